$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new column before column V (22nd column, "style list") to make room
# for the new "server_calculation" column.
$ws.Columns("V").Insert()

# Set header value + bold style matching neighboring header cells
$ws.Range("V1").Value = "server_calculation"
$ws.Range("V1").Font.Bold = $true

# Update selection / frozen pane view to mirror authored state
$ws.Range("V2").Select()
